## ETL_Concepts.xlsx - add ETL table rows for vitals, labs/imaging concepts
## (commit: "add etl tables for vitals, labs/imaging")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# ---------------------------------------------------------------------------
# New data rows appended right after the existing last row (145) of the
# Concepts table. Columns: A=concept_id, B=uuid, C=openmrs concept name.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 146; ConceptId = 5089; Uuid = "a8a660ca-1350-11df-a1f1-0026b9348838"; Name = "WEIGHT (KG)" },
    @{ Row = 147; ConceptId = 5090; Uuid = "a8a6619c-1350-11df-a1f1-0026b9348838"; Name = "HEIGHT (CM)" },
    @{ Row = 148; ConceptId = 5088; Uuid = "a8a65fee-1350-11df-a1f1-0026b9348838"; Name = "TEMPERATURE (C)" },
    @{ Row = 149; ConceptId = 5092; Uuid = "a8a66354-1350-11df-a1f1-0026b9348838"; Name = "BLOOD OXYGEN SATURATION" },
    @{ Row = 150; ConceptId = 5085; Uuid = "a8a65d5a-1350-11df-a1f1-0026b9348838"; Name = "SYSTOLIC BLOOD PRESSURE" },
    @{ Row = 151; ConceptId = 5086; Uuid = "a8a65e36-1350-11df-a1f1-0026b9348838"; Name = "DIASTOLIC BLOOD PRESSURE" },
    @{ Row = 152; ConceptId = 5087; Uuid = "a8a65f12-1350-11df-a1f1-0026b9348838"; Name = "PULSE" },
    @{ Row = 153; ConceptId = 21;   Uuid = "a8908a16-1350-11df-a1f1-0026b9348838"; Name = "HEMOGLOBIN" },
    @{ Row = 154; ConceptId = 653;  Uuid = "a896c8ae-1350-11df-a1f1-0026b9348838"; Name = "AST" },
    @{ Row = 155; ConceptId = 790;  Uuid = "a897e450-1350-11df-a1f1-0026b9348838"; Name = "SERUM CREATININE" },
    @{ Row = 156; ConceptId = 12;   Uuid = "a8908192-1350-11df-a1f1-0026b9348838"; Name = "X-RAY, CHEST, PRELIMINARY FINDINGS" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.ConceptId
    $ws.Cells.Item($r.Row, 2).Value = $r.Uuid
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
}

# Match the look of the rows above (font/size) by copying the formatting
# from the last pre-existing data row (145) down onto the newly written
# rows, so the new cells pick up the same styles as the rest of the table
# (instead of creating brand new style/font entries).
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$ws.Range("A145").Copy()
$ws.Range("A146:A156").PasteSpecial($xlPasteFormats)

$ws.Range("B145:C145").Copy()
$ws.Range("B146:C156").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# View/selection state (scroll position + active cell), as recorded when the
# workbook was last saved.
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 54
$ws.Range("C63").Select()
